$wb = $excel.ActiveWorkbook

# Row -> new F value (想去人数) for both the "展览" and "全部类型" sheets.
# Row 13's new value differs slightly between the two sheets.
$fUpdatesCommon = @{
    3  = 11371
    4  = 10706
    5  = 599
    7  = 770
    8  = 113
    9  = 43
    11 = 32
    12 = 10544
    17 = 23
    18 = 98
    19 = 412
    20 = 11088
    21 = 10837
    23 = 21
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # G2: 250 (number) -> "不可售" (text)
    $ws.Range("G2").Value = "不可售"

    foreach ($row in $fUpdatesCommon.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdatesCommon[$row]
    }
}

# Row 13's F value diverges between the two sheets (3344 vs 3345).
$wb.Worksheets.Item("展览").Cells.Item(13, 6).Value = 3344
$wb.Worksheets.Item("全部类型").Cells.Item(13, 6).Value = 3345
